# Auto-generated edit script: update Sheets via scheduled runner
# Applies updated currentAveragePrice / Leve market values across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 537.7143
$ws.Range("I2").Value = 405
$ws.Range("K2").Value = 405
$ws.Range("M2").Value = -292

$ws.Range("H33").Value = 143.7
$ws.Range("I33").Value = 97.833336
$ws.Range("K33").Value = 97.833336
$ws.Range("M33").Value = 131.166664

$ws.Range("H49").Value = 1900
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws.Range("H53").Value = 1327.1666
$ws.Range("I53").Value = 636.2222
$ws.Range("J53").Value = 3400
$ws.Range("K53").Value = 636.2222
$ws.Range("L53").Value = 3400
$ws.Range("M53").Value = 0.7777999999999565
$ws.Range("N53").Value = -4674

$ws.Range("H99").Value = 1290.8334
$ws.Range("I99").Value = 379.83334
$ws.Range("J99").Value = 2201.8333
$ws.Range("K99").Value = 1139.50002
$ws.Range("L99").Value = 6605.499899999999
$ws.Range("M99").Value = 358.4999800000001
$ws.Range("N99").Value = -9601.499899999999

$ws.Range("H131").Value = 3292.8
$ws.Range("I131").Value = 3379.125
$ws.Range("K131").Value = 10137.375
$ws.Range("M131").Value = -5097.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2141.6667
$ws.Range("I45").Value = 1313
$ws.Range("K45").Value = 1313
$ws.Range("M45").Value = -936

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1350
$ws.Range("I5").Value = 1350
$ws.Range("K5").Value = 1350
$ws.Range("M5").Value = -1237

$ws.Range("H22").Value = 455.2857
$ws.Range("I22").Value = 447.25
$ws.Range("J22").Value = 466
$ws.Range("K22").Value = 447.25
$ws.Range("L22").Value = 466
$ws.Range("M22").Value = -274.25
$ws.Range("N22").Value = -812

$ws.Range("H86").Value = 4387.3335
$ws.Range("I86").Value = 1820.8572
$ws.Range("J86").Value = 7980.4
$ws.Range("K86").Value = 1820.8572
$ws.Range("L86").Value = 7980.4
$ws.Range("M86").Value = -697.8571999999999
$ws.Range("N86").Value = -10226.4

$ws.Range("H89").Value = 4387.3335
$ws.Range("I89").Value = 1820.8572
$ws.Range("J89").Value = 7980.4
$ws.Range("K89").Value = 9104.286
$ws.Range("L89").Value = 39902
$ws.Range("M89").Value = -3488.286
$ws.Range("N89").Value = -51134

$ws.Range("H105").Value = 2231.5
$ws.Range("I105").Value = 2258.2
$ws.Range("K105").Value = 2258.2
$ws.Range("M105").Value = -511.1999999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 7443.125
$ws.Range("I6").Value = 363.2857
$ws.Range("J6").Value = 57002
$ws.Range("K6").Value = 363.2857
$ws.Range("L6").Value = 57002
$ws.Range("M6").Value = -250.2857
$ws.Range("N6").Value = -57228

$ws.Range("H59").Value = 34373.8
$ws.Range("I59").Value = 12552
$ws.Range("K59").Value = 12552
$ws.Range("M59").Value = -11407

$ws.Range("H68").Value = 71941.25
$ws.Range("J68").Value = 71941.25
$ws.Range("L68").Value = 71941.25
$ws.Range("N68").Value = -73439.25

$ws.Range("H71").Value = 71941.25
$ws.Range("J71").Value = 71941.25
$ws.Range("L71").Value = 215823.75
$ws.Range("N71").Value = -223311.75

$ws.Range("H95").Value = 10575
$ws.Range("J95").Value = 10575
$ws.Range("L95").Value = 10575
$ws.Range("N95").Value = -16067

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 339.66666
$ws.Range("J86").Value = 779
$ws.Range("L86").Value = 2337
$ws.Range("N86").Value = -4709

$ws.Range("H89").Value = 339.66666
$ws.Range("J89").Value = 779
$ws.Range("L89").Value = 7011
$ws.Range("N89").Value = -18867

$ws.Range("H107").Value = 481.05554
$ws.Range("I107").Value = 401.5
$ws.Range("K107").Value = 1204.5
$ws.Range("M107").Value = 715.5

$ws.Range("H113").Value = 1522.3334
$ws.Range("I113").Value = 947.75
$ws.Range("K113").Value = 2843.25
$ws.Range("M113").Value = -673.25

$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws.Range("H138").Value = 5156.5557
$ws.Range("I138").Value = 1264.75
$ws.Range("J138").Value = 8270
$ws.Range("K138").Value = 3794.25
$ws.Range("L138").Value = 24810
$ws.Range("M138").Value = 1345.75
$ws.Range("N138").Value = -35090

$ws.Range("H139").Value = 1610.5294
$ws.Range("I139").Value = 957.4286
$ws.Range("K139").Value = 2872.2858
$ws.Range("M139").Value = 2267.7142

$ws.Range("H140").Value = 2803.5264
$ws.Range("I140").Value = 2515.7058
$ws.Range("K140").Value = 7547.117400000001
$ws.Range("M140").Value = -2367.117400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H97").Value = 427.9
$ws.Range("I97").Value = 319.75
$ws.Range("K97").Value = 319.75
$ws.Range("M97").Value = 176.25

$ws.Range("H101").Value = 25374.5
$ws.Range("J101").Value = 25374.5
$ws.Range("L101").Value = 25374.5
$ws.Range("N101").Value = -31864.5

$ws.Range("H109").Value = 34567
$ws.Range("J109").Value = 34567
$ws.Range("L109").Value = 34567
$ws.Range("N109").Value = -36647

$ws.Range("H122").Value = 1899.6875
$ws.Range("J122").Value = 1597
$ws.Range("L122").Value = 4791
$ws.Range("N122").Value = -9691

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 770
$ws.Range("I9").Value = 770
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 770
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -546
$ws.Range("N9").ClearContents()

$ws.Range("H46").Value = 2543.9583
$ws.Range("J46").Value = 2702.5264
$ws.Range("L46").Value = 2702.5264
$ws.Range("N46").Value = -3078.5264

$ws.Range("H61").Value = 4037.6667
$ws.Range("I61").Value = 1530.6
$ws.Range("K61").Value = 1530.6
$ws.Range("M61").Value = -1328.6

$ws.Range("H113").Value = 4037.6667
$ws.Range("I113").Value = 1530.6
$ws.Range("K113").Value = 1530.6
$ws.Range("M113").Value = 639.4000000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
